$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ======================================================================
# Step 1: create the new cells (B17/C17, B22/C22, A23, B25/C25) by first
# copying formatting from a neighbouring cell in the same column, so they
# pick up the correct style (bold label / wrapped text / red wrapped text)
# instead of the sheet default.
# ======================================================================

$ws.Range("B16").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("C16").Copy()
$ws.Range("C17").PasteSpecial(-4122)

$ws.Range("B21").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("C21").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("A22").Copy()
$ws.Range("A23").PasteSpecial(-4122)

$ws.Range("B24").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$ws.Range("C24").Copy()
$ws.Range("C25").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ======================================================================
# Step 2: write the final text for every row, from the bottom up so that
# no target value is overwritten by a later step, and content that needs
# to move to a new row is copied from the old row before it is cleared.
# ======================================================================

# Row 25 (new): requisito 2 (was on row 24)
$ws.Range("B25").Value = "LOT2053 -  Microbiologia  (Requisito fraco)`n"
$ws.Range("C25").Value = "LOT2053 -  Microbiologia  (Requisito fraco)`n"

# Row 24: requisito 1 (was on row 23)
$ws.Range("B24").Value = "LOT2008 -  Bioquímica II  (Requisito fraco)`n"
$ws.Range("C24").Value = "LOT2008 -  Bioquímica II  (Requisito fraco)`n"

# Row 23 (new): "Requisitos:" label moves down from row 22; clear old B23/C23
$ws.Range("A23").Value = "Requisitos:"
$ws.Range("B23").Clear()
$ws.Range("C23").Clear()

# Row 22: "Bibliografia:" label moves down from row 21, with new bibliography text
$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = "1. AMERINE, M.A, OUGH,C.S., Methods for analysis of musts and wines. New York: John Wiley & Sons, 1980. `n2. AMORIM, H.V., Fermentação Alcoólica ciência e tecnologia. Piracicaba: Fermentec, 2006.`n3. BORZANI, W., SCHMIDELL, W., LIMA, U.A., AQUARONE, E. Série de Biotecnologia Vol. 1  Fundamentos e Vol. 4 Processos Fermentativos e Enzimáticos. São Paulo: Ed.Edgard Blucher, 2001.`n4. EL-MANSI, E.M.T., BRYCE, C.E.A., DEMAIN, A.L., ALLMAN,A.R. Fermentation Microbiology and Biotechnology. 2ª Ed. New York: CRC Taylor & Francis, 2007."
$ws.Range("C22").Value = "1. AMERINE, M.A, OUGH,C.S., Methods for analysis of musts and wines. New York: John Wiley & Sons, 1980. `n2. AMORIM, H.V., Fermentação Alcoólica ciência e tecnologia. Piracicaba: Fermentec, 2006.`n3. BORZANI, W., SCHMIDELL, W., LIMA, U.A., AQUARONE, E. Série de Biotecnologia Vol. 1  Fundamentos e Vol. 4 Processos Fermentativos e Enzimáticos. São Paulo: Ed.Edgard Blucher, 2001.`n4. EL-MANSI, E.M.T., BRYCE, C.E.A., DEMAIN, A.L., ALLMAN,A.R. Fermentation Microbiology and Biotechnology. 2ª Ed. New York: CRC Taylor & Francis, 2007."

# Row 21: "Norma de recuperação:" label moves down from row 20, formula text moves up from row 21
$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Range("B21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Range("C21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"

# Row 20: "Critério:" label moves down from row 19, formula text moves up from row 20
$ws.Range("A20").Value = "Critério:"
$ws.Range("B20").Value = "A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + P2) / 2"
$ws.Range("C20").Value = "A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + P2) / 2"

# Row 19: "Método:" label moves down from row 18, text moves up from row 19
$ws.Range("A19").Value = "Método:"
$ws.Range("B19").Value = "A avaliação será feita por meio de provas escritas."
$ws.Range("C19").Value = "A avaliação será feita por meio de provas escritas."

# Row 18: "Avaliação:" label moves down from row 17; clear old B18/C18 (stray duplicate text)
$ws.Range("A18").Value = "Avaliação:"
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()

# Row 17: "Syllabus:" label moves down from row 16; English full syllabus moves up from row 16 into new B17/C17
$ws.Range("A17").Value = "Syllabus:"
$ws.Range("B17").Value = "1.Biotechnology: concepts, application areas, multidisciplinary characteristic and examples of biotechnological products and processes.2.Fermentative processes: concept, enzymatic and fermentative processes, steps of fermentative process (downstream x upstream). Fermentative process modalities: a) batch and fed-batch fermentation, semi continuous and, continuous processes; b) induced and spontaneous fermentation;  c) semi solid fermentation;  d) oxygen supply; e) submerged and in surface processes; f) kinetics of the product formation in relation to the primary metabolism according to Gaden.3.Biochemistry of the fermentation: Fermentation – concepts, objectives, aerobic x anaerobic metabolisms; energy balance; preliminary steps of fermentation (extracellular hydrolysis and membrane permeability); metabolic pathways of industrial interest: a) EMP pathway; reactions and  allosteric control; alcoholic fermentation, homolactic fermentation, acetone/butanol, mixed-acid and 2,3 butanediol; b) Fosfo-Ketolase pathway; heterolactic fermentation and c) Entner Doudoroff pathway: alcoholic fermentation by Zymmonas mobilis. Fermentation balance: % recovered carbon and oxi-redox balance; Evaluation parameters of a fermentative process: yield, fermentation efficiency and productivity. Processes of interest: cocoa processing, ethanol production, fermented food and others."
$ws.Range("C17").Value = "1.Biotechnology: concepts, application areas, multidisciplinary characteristic and examples of biotechnological products and processes.2.Fermentative processes: concept, enzymatic and fermentative processes, steps of fermentative process (downstream x upstream). Fermentative process modalities: a) batch and fed-batch fermentation, semi continuous and, continuous processes; b) induced and spontaneous fermentation;  c) semi solid fermentation;  d) oxygen supply; e) submerged and in surface processes; f) kinetics of the product formation in relation to the primary metabolism according to Gaden.3.Biochemistry of the fermentation: Fermentation – concepts, objectives, aerobic x anaerobic metabolisms; energy balance; preliminary steps of fermentation (extracellular hydrolysis and membrane permeability); metabolic pathways of industrial interest: a) EMP pathway; reactions and  allosteric control; alcoholic fermentation, homolactic fermentation, acetone/butanol, mixed-acid and 2,3 butanediol; b) Fosfo-Ketolase pathway; heterolactic fermentation and c) Entner Doudoroff pathway: alcoholic fermentation by Zymmonas mobilis. Fermentation balance: % recovered carbon and oxi-redox balance; Evaluation parameters of a fermentative process: yield, fermentation efficiency and productivity. Processes of interest: cocoa processing, ethanol production, fermented food and others."

# Row 16: "Programa:" label moves down from row 15; new Portuguese full syllabus goes into B16/C16
$ws.Range("A16").Value = "Programa:"
$ws.Range("B16").Value = "1. Biotecnologia: conceitos, áreas de aplicação, caráter multidisciplinar e exemplos de produtos biotecnológicos.2. Processos Fermentativos: conceito, exemplos, fases de um processo fermentativo. Modalidades de Processos Fermentativos: a)formas de condução; b) fermentação induzida e espontânea; c) estado físico do meio de fermentação; d) suprimento de oxigênio; e) processos submersos e em superfície; f) cinética de formação de produto em relação do metabolismo primário.3. Bioquímica das fermentações: fermentação  conceitos, objetivos, aerobiose x anaerobiose; balanço energético; estágios preliminares da fermentação (hidrólise extracelular e permeabilidade da membrana); vias metabólica de interesse industrial: a) via glicolítica: reações e controle; fermentação alcoólica, homoláctica, acetona/butanol, ácido-mista e 2,3 butanodiol; b) Via Fosfo-Cetolase: fermentação heteroláctica e c) via Entner Doudoroff: fermentação alcoólica por Zymomonas mobilis. Balanço da Fermentação: % de carbono recuperado e balanço de oxi-redução; parâmetros de avaliação - rendimento, eficiência e produtividade de processos fermentativos; Processos de Interesse: processamento de cacau, produção de etanol, alimentos fermentados e outros."
$ws.Range("C16").Value = "1. Biotecnologia: conceitos, áreas de aplicação, caráter multidisciplinar e exemplos de produtos biotecnológicos.2. Processos Fermentativos: conceito, exemplos, fases de um processo fermentativo. Modalidades de Processos Fermentativos: a)formas de condução; b) fermentação induzida e espontânea; c) estado físico do meio de fermentação; d) suprimento de oxigênio; e) processos submersos e em superfície; f) cinética de formação de produto em relação do metabolismo primário.3. Bioquímica das fermentações: fermentação  conceitos, objetivos, aerobiose x anaerobiose; balanço energético; estágios preliminares da fermentação (hidrólise extracelular e permeabilidade da membrana); vias metabólica de interesse industrial: a) via glicolítica: reações e controle; fermentação alcoólica, homoláctica, acetona/butanol, ácido-mista e 2,3 butanodiol; b) Via Fosfo-Cetolase: fermentação heteroláctica e c) via Entner Doudoroff: fermentação alcoólica por Zymomonas mobilis. Balanço da Fermentação: % de carbono recuperado e balanço de oxi-redução; parâmetros de avaliação - rendimento, eficiência e produtividade de processos fermentativos; Processos de Interesse: processamento de cacau, produção de etanol, alimentos fermentados e outros."

# Row 15: "Short syllabus:" label moves down from row 14; English short syllabus stays (moved up from row 14)
$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("B15").Value = "Biotechnology (field of applications); fermentative processes; biochemistry of the fermentations (metabolic pathways of industrial interest); fermentative processes of industrial interest"
$ws.Range("C15").Value = "Biotechnology (field of applications); fermentative processes; biochemistry of the fermentations (metabolic pathways of industrial interest); fermentative processes of industrial interest"

# Row 14: "Programa resumido:" label moves down from row 13; new Portuguese short syllabus goes into B14/C14
$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = "Biotecnologia; processos fermentativos; bioquímica das fermentações (vias metabólicas de interesse industrial); processos fermentativos de interesse industrial."
$ws.Range("C14").Value = "Biotecnologia; processos fermentativos; bioquímica das fermentações (vias metabólicas de interesse industrial); processos fermentativos de interesse industrial."

# Row 13: old "Programa resumido:" label removed, teacher record moves up into B13/C13
$ws.Range("A13").Clear()
$ws.Range("B13").Value = "3403572 - Ismael Maciel de Mancilha"
$ws.Range("C13").Value = "3403572 - Ismael Maciel de Mancilha"

# Row 10: "Objetivos:" gets the new Portuguese objectives text
$ws.Range("B10").Value = "Levar aos estudantes conhecimentos básicos sobre: a) processos fermentativos, com ênfase em processos de interesse industrial; b) bioquímica das fermentações focando as rotas metabólicas utilizadas por microrganismos de interesse industrial; c) suas respectivas aplicações em processos industriais, permitindo a determinação de parâmetros de avaliação de desempenho."
$ws.Range("C10").Value = "Levar aos estudantes conhecimentos básicos sobre: a) processos fermentativos, com ênfase em processos de interesse industrial; b) bioquímica das fermentações focando as rotas metabólicas utilizadas por microrganismos de interesse industrial; c) suas respectivas aplicações em processos industriais, permitindo a determinação de parâmetros de avaliação de desempenho."

# ======================================================================
# Step 3: fix up row heights to match the new layout.
# ======================================================================
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).AutoFit()
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 120
$ws.Rows.Item(23).AutoFit()
$ws.Rows.Item(25).RowHeight = 30

